$wb = $excel.ActiveWorkbook

# ALC row 33 (diff @ 2219)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 306.70587
$ws.Range("I33").Value = 322.75
$ws.Range("J33").Value = 50
$ws.Range("K33").Value = 322.75
$ws.Range("L33").Value = 50
$ws.Range("M33").Value = -93.75
$ws.Range("N33").Value = -508

# ALC row 41 (diff @ 2614)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 433.33334
$ws.Range("I41").Value = 433.33334
$ws.Range("K41").Value = 433.33334
$ws.Range("M41").Value = 6.666659999999979

# ALC row 98 (diff @ 5392)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 732.6842
$ws.Range("I98").Value = 701.2353000000001
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 701.2353000000001
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = 796.7646999999999
$ws.Range("N98").Value = -3996

# ALC row 99 (diff @ 5444)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 1500
$ws.Range("I99").Value = 1500
$ws.Range("K99").Value = 4500
$ws.Range("M99").Value = -3002

# ALC row 100 (diff @ 5493)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3335327.2
$ws.Range("I100").Value = 3335327.2
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 3335327.2
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -3334786.2

# ALC row 106 (diff @ 5790)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 7249.4
$ws.Range("J106").Value = 7999.25
$ws.Range("L106").Value = 7999.25
$ws.Range("N106").Value = -9261.25

# ALC row 122 (diff @ 6580)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 732.6842
$ws.Range("I122").Value = 701.2353000000001
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2103.7059
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = 346.2941000000001
$ws.Range("N122").Value = -7900

# ALC row 132 (diff @ 7082)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 812.7907
$ws.Range("I132").Value = 808.3333
$ws.Range("K132").Value = 2424.9999
$ws.Range("M132").Value = 105.0001000000002

# ALC row 137 (diff @ 7330)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1799.7646
$ws.Range("I137").Value = 1437.3914
$ws.Range("K137").Value = 4312.174199999999
$ws.Range("M137").Value = -1762.174199999999

# ALC row 138 (diff @ 7382)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4286.853
$ws.Range("J138").Value = 4358.567
$ws.Range("L138").Value = 13075.701
$ws.Range("N138").Value = -23355.701

# ARM row 32 (diff @ 9112)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4896.891
$ws.Range("I32").Value = 4554.284
$ws.Range("J32").Value = 12434.25
$ws.Range("K32").Value = 4554.284
$ws.Range("L32").Value = 12434.25
$ws.Range("M32").Value = -4267.284
$ws.Range("N32").Value = -13008.25

# ARM row 45 (diff @ 9731)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2746.5
$ws.Range("I45").Value = 2746.5
$ws.Range("K45").Value = 2746.5
$ws.Range("M45").Value = -2369.5

# ARM row 97 (diff @ 12177)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 254.4
$ws.Range("I97").Value = 254.4
$ws.Range("K97").Value = 254.4
$ws.Range("M97").Value = 241.6

# ARM row 122 (diff @ 13357)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5885.5713
$ws.Range("I122").Value = 6899.8
$ws.Range("J122").Value = 3350
$ws.Range("K122").Value = 20699.4
$ws.Range("L122").Value = 10050
$ws.Range("M122").Value = -18249.4
$ws.Range("N122").Value = -14950

# ARM row 132 (diff @ 13838)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2743.6086
$ws.Range("I132").Value = 1940.5333
$ws.Range("K132").Value = 5821.5999
$ws.Range("M132").Value = -3291.5999

# BSM row 20 (diff @ 15277)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 995
$ws.Range("I20").Value = 995
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 995
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -748

# BSM row 99 (diff @ 19040)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 937
$ws.Range("I99").Value = 950
$ws.Range("K99").Value = 950
$ws.Range("M99").Value = 548

# BSM row 134 (diff @ 20683)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2896.8
$ws.Range("I134").Value = 2969.2
$ws.Range("J134").Value = 2679.6
$ws.Range("K134").Value = 8907.599999999999
$ws.Range("L134").Value = 8038.799999999999
$ws.Range("M134").Value = -6372.599999999999
$ws.Range("N134").Value = -13108.8

# CRP row 86 (diff @ 25183)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 22454.334
$ws.Range("J86").Value = 31755.4
$ws.Range("L86").Value = 31755.4
$ws.Range("N86").Value = -34001.4

# CRP row 89 (diff @ 25330)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 22454.334
$ws.Range("J89").Value = 31755.4
$ws.Range("L89").Value = 158777
$ws.Range("N89").Value = -170009

# CUL row 80 (diff @ 31810)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 205199.4
$ws.Range("I80").Value = 4999.5
$ws.Range("J80").Value = 338666
$ws.Range("K80").Value = 14998.5
$ws.Range("L80").Value = 1015998
$ws.Range("M80").Value = -14062.5
$ws.Range("N80").Value = -1017870

# CUL row 83 (diff @ 31954)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 205199.4
$ws.Range("I83").Value = 4999.5
$ws.Range("J83").Value = 338666
$ws.Range("K83").Value = 44995.5
$ws.Range("L83").Value = 3047994
$ws.Range("M83").Value = -40315.5
$ws.Range("N83").Value = -3057354

# CUL row 86 (diff @ 32098)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

# CUL row 89 (diff @ 32242)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

# CUL row 122 (diff @ 33853)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 60416.53
$ws.Range("I122").Value = 1498
$ws.Range("J122").Value = 78545.30499999999
$ws.Range("K122").Value = 13482
$ws.Range("L122").Value = 706907.7449999999
$ws.Range("M122").Value = -11032
$ws.Range("N122").Value = -711807.7449999999

# CUL row 131 (diff @ 34300)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1530.6666
$ws.Range("J131").Value = 1999
$ws.Range("L131").Value = 5997
$ws.Range("N131").Value = -16077

# CUL row 140 (diff @ 34756)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1066.1333
$ws.Range("I140").Value = 1066.1333
$ws.Range("K140").Value = 3198.3999
$ws.Range("M140").Value = 1981.6001

# GSM row 80 (diff @ 38623)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4012.5
$ws.Range("I80").Value = 2849
$ws.Range("K80").Value = 2849
$ws.Range("M80").Value = -1851

# GSM row 83 (diff @ 38767)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4012.5
$ws.Range("I83").Value = 2849
$ws.Range("K83").Value = 14245
$ws.Range("M83").Value = -9253

# GSM row 97 (diff @ 39423)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 415.65518
$ws.Range("I97").Value = 430.9091
$ws.Range("J97").Value = 367.7143
$ws.Range("K97").Value = 430.9091
$ws.Range("L97").Value = 367.7143
$ws.Range("M97").Value = 65.09089999999998
$ws.Range("N97").Value = -1359.7143

# GSM row 102 (diff @ 39671)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2222
$ws.Range("J102").Value = 2222
$ws.Range("L102").Value = 2222
$ws.Range("N102").Value = -5466

# GSM row 122 (diff @ 40609)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1991
$ws.Range("I122").Value = 2072.8333
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 6218.499899999999
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -3768.499899999999
$ws.Range("N122").Value = -9400

# GSM row 132 (diff @ 41093)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4099.25
$ws.Range("I132").Value = 3699.3333
$ws.Range("K132").Value = 11097.9999
$ws.Range("M132").Value = -8567.999899999999

# LTW row 100 (diff @ 46287)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1331.3334
$ws.Range("I100").Value = 497
$ws.Range("K100").Value = 497
$ws.Range("M100").Value = 44

# LTW row 132 (diff @ 47810)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5081.4736
$ws.Range("I132").Value = 4918.1333
$ws.Range("K132").Value = 14754.3999
$ws.Range("M132").Value = -12224.3999

# LTW row 136 (diff @ 48003)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3477.7354
$ws.Range("I136").Value = 3723.7407
$ws.Range("J136").Value = 2528.8572
$ws.Range("K136").Value = 11171.2221
$ws.Range("L136").Value = 7586.571599999999
$ws.Range("M136").Value = -8621.222099999999
$ws.Range("N136").Value = -12686.5716

# WVR row 81 (diff @ 52124)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4399.7646
$ws.Range("I81").Value = 2621.5715
$ws.Range("K81").Value = 5243.143
$ws.Range("M81").Value = -4182.143

# WVR row 84 (diff @ 52268)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 4399.7646
$ws.Range("I84").Value = 2621.5715
$ws.Range("K84").Value = 26215.715
$ws.Range("M84").Value = -20911.715

# WVR row 132 (diff @ 54542)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3427.158
$ws.Range("I132").Value = 2670.4443
$ws.Range("J132").Value = 4108.2
$ws.Range("K132").Value = 8011.3329
$ws.Range("L132").Value = 12324.6
$ws.Range("M132").Value = -5481.3329
$ws.Range("N132").Value = -17384.6
